$d = $word.ActiveDocument
$nl = [char]10

# Table 2 (1-indexed) = מתמטיקה (Math)
$mathComment = "במחצית למדנו משוואות ב2 נעלמים, פיתחנו כמה שיטות לבעיה זו,בנוסף התעסקנו בבעיות תנועה וזמן ולמדנו איך להתמודד מול זאת" + $nl + "היה לנו הספקים מעולים!" + $nl + "גל את ילדה מקסימה, שיהיה לך הרבה הצלחה בהמשך! "
$d.Tables(2).Cell(1, 2).Range.Text = $mathComment
$d.Tables(2).Cell(2, 1).Range.Paragraphs(2).Range.Text = "92"

# Table 3 (1-indexed) = אנגלית (English)
$englishComment = "במחצית זאת התמקדנו על הבנה חזקה של הטקסטים ולמדנו את השיטות להבנת הנקרא, חזרנו על שאלות חוזרות ופיתחנו שיטות קלות לפיתרתן." + $nl + "גל את ילדה נהדרת, הרבה הצלחה!"
$d.Tables(3).Cell(1, 2).Range.Text = $englishComment
$d.Tables(3).Cell(2, 1).Range.Paragraphs(2).Range.Text = "87"

# Table 4 (1-indexed) = תולדות ישראל (Jewish History)
$historyComment = "במחצית זאת למדנו על גדולי ישראל בכל מיני יבשות, על המצב של היהודים בתקופות שלטון שונות," + $nl + "גל הרבה הצלחה!"
$d.Tables(4).Cell(1, 2).Range.Text = $historyComment
$d.Tables(4).Cell(2, 1).Range.Paragraphs(2).Range.Text = "91"
